# Fixes duplicate-match data entry error.
# Re-maps the India (A:F) and New Zealand (J:O) scorecards to the correct
# batting order / runs / balls / dismissal / bowler figures, and recomputes the
# totals and bowling-figures tables underneath (rows 16 and 21-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validation")

# --- row 2: Batter 1 ---
$ws.Range("A2").Value = "KL Rahul"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "Caught"
$ws.Range("E2").Value = " Ish Sodhi"
$ws.Range("J2").Value = "Martin Guptill"
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 7
$ws.Range("N2").Value = " Hardik Pandya"
# --- row 3: Batter 2 ---
$ws.Range("A3").Value = "Rohit Sharma"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "Bowled"
$ws.Range("E3").Value = " Trent Boult"
$ws.Range("J3").Value = "Daryl Mitchell"
$ws.Range("K3").Value = 7
$ws.Range("L3").Value = 7
$ws.Range("M3").Value = "Bowled"
$ws.Range("N3").Value = " Jasprit Bumrah"
# --- row 4: Batter 3 ---
$ws.Range("A4").Value = "Virat Kohli(C)"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = " Tim Southee"
$ws.Range("J4").Value = "Kane Williamson(C)"
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 2
$ws.Range("N4").Value = " Jasprit Bumrah"
# --- row 5: Batter 4 ---
$ws.Range("A5").Value = "Suryakumar Yadav"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "LBW"
$ws.Range("E5").Value = " Tim Southee"
$ws.Range("J5").Value = "Devon Conway"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = "Caught"
$ws.Range("N5").Value = " Jasprit Bumrah"
# --- row 6: Batter 5 ---
$ws.Range("A6").Value = "Rishabh Pant"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("E6").Value = " Tim Southee"
$ws.Range("J6").Value = "Glenn Phillips"
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 2
$ws.Range("N6").Value = " Mohommad Shami"
# --- row 7: Batter 6 ---
$ws.Range("A7").Value = "Ravindra Jadeja"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "NOT OUT"
$ws.Range("E7").Value = " "
$ws.Range("J7").Value = "James Neesham"
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = "Bowled"
$ws.Range("N7").Value = " Mohommad Shami"
# --- row 8: Batter 7 ---
$ws.Range("A8").Value = "Hardik Pandya"
$ws.Range("B8").Value = 11
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "LBW"
$ws.Range("E8").Value = " Ish Sodhi"
$ws.Range("J8").Value = "Mitchell Santner"
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 4
$ws.Range("N8").Value = " Bhuvneshwar Kumar"
# --- row 9: Batter 8 ---
$ws.Range("A9").Value = "Bhuvneshwar Kumar"
$ws.Range("B9").Value = 13
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "LBW"
$ws.Range("E9").Value = " Adam Milne"
$ws.Range("J9").Value = "Adam Milne"
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1
$ws.Range("N9").Value = " Bhuvneshwar Kumar"
# --- row 10: Batter 9 ---
$ws.Range("A10").Value = "Mohommad Shami"
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "Bowled"
$ws.Range("E10").Value = " Adam Milne"
$ws.Range("J10").Value = "Ish Sodhi"
$ws.Range("K10").Value = 4
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = "Caught"
$ws.Range("N10").Value = " Jasprit Bumrah"
# --- row 11: Batter 10 ---
$ws.Range("A11").Value = "Jasprit Bumrah"
$ws.Range("B11").Value = 24
$ws.Range("C11").Value = 7
$ws.Range("E11").Value = " Mitchell Santner"
$ws.Range("J11").Value = "Tim Southee"
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1
$ws.Range("N11").Value = " Hardik Pandya"
# --- row 12: Batter 11 ---
$ws.Range("A12").Value = "Yuzvendra Chahal"
$ws.Range("B12").Value = 0
$ws.Range("D12").Value = "LBW"
$ws.Range("E12").Value = " Trent Boult"
$ws.Range("J12").Value = "Trent Boult"
$ws.Range("K12").Value = 23
$ws.Range("L12").Value = 7
$ws.Range("M12").Value = "NOT OUT"
$ws.Range("N12").Value = " "
# --- row 16: Innings totals ---
$ws.Range("A16").Value = 88
$ws.Range("C16").Value = "'5.1"
$ws.Range("D16").Value = 31
$ws.Range("J16").Value = 66
$ws.Range("L16").Value = "'6.1"
$ws.Range("M16").Value = 37
# --- row 21: Bowler 1 ---
$ws.Range("A21").Value = "Tim Southee"
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 12
$ws.Range("J21").Value = "Mohommad Shami"
$ws.Range("L21").Value = 15
$ws.Range("M21").Value = 2
$ws.Range("N21").Value = 15
# --- row 22: Bowler 2 ---
$ws.Range("A22").Value = "Ish Sodhi"
$ws.Range("C22").Value = 18
$ws.Range("E22").Value = 18
$ws.Range("J22").Value = "Bhuvneshwar Kumar"
$ws.Range("K22").Value = "'1.0"
$ws.Range("L22").Value = 9
$ws.Range("M22").Value = 2
$ws.Range("N22").Value = 9
# --- row 23: Bowler 3 ---
$ws.Range("A23").Value = "Adam Milne"
$ws.Range("C23").Value = 14
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 14
$ws.Range("J23").Value = "Hardik Pandya"
$ws.Range("K23").Value = "'1.0"
$ws.Range("L23").Value = 15
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 15
# --- row 24: Bowler 4 ---
$ws.Range("A24").Value = "Mitchell Santner"
$ws.Range("B24").Value = "'1.0"
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 22
$ws.Range("J24").Value = "Yuzvendra Chahal"
$ws.Range("L24").Value = 19
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 9.5
# --- row 25: Bowler 5 ---
$ws.Range("A25").Value = "Trent Boult"
$ws.Range("B25").Value = "'1.1"
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 20
$ws.Range("J25").Value = "Jasprit Bumrah"
$ws.Range("K25").Value = "'1.1"
$ws.Range("L25").Value = 8
$ws.Range("N25").Value = 7.27
